$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from 45204 to
# 45205 for every data row (rows 2 through 103).
$range = $ws.Range("C2:C103")
$range.Value2 = 45205
